$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need the Text number format
# forced first, so Excel stores them as the exact literal string (e.g. "1.00")
# instead of silently converting them to a numeric value (e.g. 1).

$ws.Range("D2").Value = "64.757.33"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "2.531.02"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.41"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.55"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "2.532.97"
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.33"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000181"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "2.990.07"
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "64.317.16"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").Value = "2.537.92"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.93"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.99"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.28"
$ws.Range("E21").Value = "  +3.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "330.28"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.24"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.04"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.66"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "649.04"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  +7.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.51"
$ws.Range("E30").Value = "  +5.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.56"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.85"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.58"
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "155.08"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.99"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.373"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("E43").Value = "  +3.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "163.67"
$ws.Range("E44").Value = "  +7.38%  "
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.61"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.65"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.56"
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.622"
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("E51").Value = "  +1.30%  "
